# Weekly update: insert a new price-report row at row 22 (most-recent-first
# dataset), shifting all existing data rows down by one. The new row mostly
# mirrors the (now-shifted) row below it, with a handful of fields updated
# to reflect the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 22; Excel shifts rows 22:113 down to 23:114
# and extends the used range / dimension to A1:R114 automatically.
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with this week's record.
$ws.Range("A22").Value2 = 7
$ws.Range("B22").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value2 = "Ñuble"
$ws.Range("D22").Value2 = 45145
$ws.Range("E22").Value2 = 16
$ws.Range("F22").Value2 = 100112013
$ws.Range("G22").Value2 = "Alcachofa"
$ws.Range("H22").Value2 = "Madrigal"
$ws.Range("I22").Value2 = "Primera"
$ws.Range("J22").Value2 = 80
$ws.Range("K22").Value2 = 12000
$ws.Range("L22").Value2 = 12000
$ws.Range("M22").Value2 = 12000
$ws.Range("N22").Value2 = "`$/caja 50 unidades"
$ws.Range("O22").Value2 = "Provincia de Limarí"
$ws.Range("P22").Value2 = 240
$ws.Range("Q22").Value2 = 50
$ws.Range("R22").Value2 = "Hortaliza"
